# Add VAT breakdown columns (Net, VAT, Gross, VAT Reclaimable) to the
# Monthly_Cashflow sheet and extend it with additional cash-flow rows
# through Jun 2026 (Q3).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monthly_Cashflow")

# --- Column widths (chars). The COM ColumnWidth property is character-
# based and Excel quantizes it to whole pixels, so values are nudged by
# -5/6 to land as close as possible to the target OOXML width.
$ws.Columns.Item(4).ColumnWidth = 11.998697916666666
$ws.Columns.Item(5).ColumnWidth = 11.998697916666666
$ws.Columns.Item(6).ColumnWidth = 14.998697916666666
$ws.Columns.Item(7).ColumnWidth = 14.998697916666666
$ws.Columns.Item(8).ColumnWidth = 29.998697916666668

# --- Cell values ---
$ws.Cells.Item(1,1).Value = "Month"
$ws.Cells.Item(1,2).Value = "Item"
$ws.Cells.Item(1,3).Value = "Net Amount"
$ws.Cells.Item(1,4).Value = "VAT (20%)"
$ws.Cells.Item(1,5).Value = "Gross Amount"
$ws.Cells.Item(1,6).Value = "VAT Reclaimable"
$ws.Cells.Item(1,7).Value = "Running Balance"
$ws.Cells.Item(1,8).Value = "Notes"

$ws.Cells.Item(2,1).Value = "Opening"
$ws.Cells.Item(2,2).Value = "Capital Raised"
$ws.Cells.Item(2,3).Value = 625000
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 625000
$ws.Cells.Item(2,6).Value = "N/A"
$ws.Cells.Item(2,7).Value = 625000
$ws.Cells.Item(2,8).Value = "Starting balance"

$ws.Cells.Item(3,1).Value = "'Dec 2025"
$ws.Cells.Item(3,2).Value = "Legal/Professional Fees"
$ws.Cells.Item(3,3).Value = 10930.5
$ws.Cells.Item(3,4).Value = 2186.1
$ws.Cells.Item(3,5).Value = 13116.6
$ws.Cells.Item(3,6).Value = "Yes"
$ws.Cells.Item(3,7).Value = 611883.4
$ws.Cells.Item(3,8).Value = "Legal and professional fees"

$ws.Cells.Item(4,1).Value = "'Dec 2025"
$ws.Cells.Item(4,2).Value = "Rent Deposit (7 months)"
$ws.Cells.Item(4,3).Value = 128388
$ws.Cells.Item(4,4).Value = 25677.6
$ws.Cells.Item(4,5).Value = 154065.6
$ws.Cells.Item(4,6).Value = "Yes"
$ws.Cells.Item(4,7).Value = 457817.8
$ws.Cells.Item(4,8).Value = "7 months deposit"

$ws.Cells.Item(5,1).Value = "'Dec 2025"
$ws.Cells.Item(5,2).Value = "Q1 Rent"
$ws.Cells.Item(5,3).Value = 32097
$ws.Cells.Item(5,4).Value = 6419.4
$ws.Cells.Item(5,5).Value = 38516.4
$ws.Cells.Item(5,6).Value = "Yes"
$ws.Cells.Item(5,7).Value = 419301.4
$ws.Cells.Item(5,8).Value = "Quarter 1 rent payment"

$ws.Cells.Item(6,1).Value = "'Dec 2025"
$ws.Cells.Item(6,2).Value = "Service Charge (Q1)"
$ws.Cells.Item(6,3).Value = 12000
$ws.Cells.Item(6,4).Value = 0
$ws.Cells.Item(6,5).Value = 12000
$ws.Cells.Item(6,6).Value = "No"
$ws.Cells.Item(6,7).Value = 407301.4
$ws.Cells.Item(6,8).Value = "Quarterly service charge - exempt"

$ws.Cells.Item(7,1).Value = "'Dec 2025"
$ws.Cells.Item(7,2).Value = "Insurance (Annual)"
$ws.Cells.Item(7,3).Value = 4800
$ws.Cells.Item(7,4).Value = 0
$ws.Cells.Item(7,5).Value = 4800
$ws.Cells.Item(7,6).Value = "No"
$ws.Cells.Item(7,7).Value = 402501.4
$ws.Cells.Item(7,8).Value = "Annual insurance - exempt"

$ws.Cells.Item(8,1).Value = "'Dec 2025"
$ws.Cells.Item(8,2).Value = "Business Rates (Dec)"
$ws.Cells.Item(8,3).Value = 5000
$ws.Cells.Item(8,4).Value = 0
$ws.Cells.Item(8,5).Value = 5000
$ws.Cells.Item(8,6).Value = "No"
$ws.Cells.Item(8,7).Value = 397501.4
$ws.Cells.Item(8,8).Value = "Monthly - no VAT on rates"

$ws.Cells.Item(9,1).Value = "'Jan 2026"
$ws.Cells.Item(9,2).Value = "Business Rates"
$ws.Cells.Item(9,3).Value = 5000
$ws.Cells.Item(9,4).Value = 0
$ws.Cells.Item(9,5).Value = 5000
$ws.Cells.Item(9,6).Value = "No"
$ws.Cells.Item(9,7).Value = 392501.4
$ws.Cells.Item(9,8).Value = "Monthly payment"

$ws.Cells.Item(10,1).Value = "'Jan 2026"
$ws.Cells.Item(10,2).Value = "Racking & Setup"
$ws.Cells.Item(10,3).Value = 20833.33
$ws.Cells.Item(10,4).Value = 4166.67
$ws.Cells.Item(10,5).Value = 25000
$ws.Cells.Item(10,6).Value = "Yes"
$ws.Cells.Item(10,7).Value = 367501.4
$ws.Cells.Item(10,8).Value = "Warehouse racking installation"

$ws.Cells.Item(11,1).Value = "'Jan 2026"
$ws.Cells.Item(11,2).Value = "Warehouse Equipment"
$ws.Cells.Item(11,3).Value = 4166.67
$ws.Cells.Item(11,4).Value = 833.33
$ws.Cells.Item(11,5).Value = 5000
$ws.Cells.Item(11,6).Value = "Yes"
$ws.Cells.Item(11,7).Value = 362501.4
$ws.Cells.Item(11,8).Value = "Initial equipment"

$ws.Cells.Item(12,1).Value = "'Feb 2026"
$ws.Cells.Item(12,2).Value = "Business Rates"
$ws.Cells.Item(12,3).Value = 5000
$ws.Cells.Item(12,4).Value = 0
$ws.Cells.Item(12,5).Value = 5000
$ws.Cells.Item(12,6).Value = "No"
$ws.Cells.Item(12,7).Value = 357501.4
$ws.Cells.Item(12,8).Value = "Monthly payment"

$ws.Cells.Item(13,1).Value = "'Mar 2026"
$ws.Cells.Item(13,2).Value = "Business Rates"
$ws.Cells.Item(13,3).Value = 5000
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 5000
$ws.Cells.Item(13,6).Value = "No"
$ws.Cells.Item(13,7).Value = 352501.4
$ws.Cells.Item(13,8).Value = "Monthly payment"

$ws.Cells.Item(14,1).Value = "'Mar 2026"
$ws.Cells.Item(14,2).Value = "Q2 Rent"
$ws.Cells.Item(14,3).Value = 32097
$ws.Cells.Item(14,4).Value = 6419.4
$ws.Cells.Item(14,5).Value = 38516.4
$ws.Cells.Item(14,6).Value = "Yes"
$ws.Cells.Item(14,7).Value = 313985
$ws.Cells.Item(14,8).Value = "Quarter 2 rent payment"

$ws.Cells.Item(15,1).Value = "'Mar 2026"
$ws.Cells.Item(15,2).Value = "Service Charge (Q2)"
$ws.Cells.Item(15,3).Value = 12000
$ws.Cells.Item(15,4).Value = 0
$ws.Cells.Item(15,5).Value = 12000
$ws.Cells.Item(15,6).Value = "No"
$ws.Cells.Item(15,7).Value = 301985
$ws.Cells.Item(15,8).Value = "Quarterly service charge"

$ws.Cells.Item(16,1).Value = "'Apr 2026"
$ws.Cells.Item(16,2).Value = "Business Rates"
$ws.Cells.Item(16,3).Value = 5000
$ws.Cells.Item(16,4).Value = 0
$ws.Cells.Item(16,5).Value = 5000
$ws.Cells.Item(16,6).Value = "No"
$ws.Cells.Item(16,7).Value = 296985
$ws.Cells.Item(16,8).Value = "Monthly payment"

$ws.Cells.Item(17,1).Value = "'May 2026"
$ws.Cells.Item(17,2).Value = "Business Rates"
$ws.Cells.Item(17,3).Value = 5000
$ws.Cells.Item(17,4).Value = 0
$ws.Cells.Item(17,5).Value = 5000
$ws.Cells.Item(17,6).Value = "No"
$ws.Cells.Item(17,7).Value = 291985
$ws.Cells.Item(17,8).Value = "Monthly payment"

$ws.Cells.Item(18,1).Value = "'Jun 2026"
$ws.Cells.Item(18,2).Value = "Business Rates"
$ws.Cells.Item(18,3).Value = 5000
$ws.Cells.Item(18,4).Value = 0
$ws.Cells.Item(18,5).Value = 5000
$ws.Cells.Item(18,6).Value = "No"
$ws.Cells.Item(18,7).Value = 286985
$ws.Cells.Item(18,8).Value = "Monthly payment"

$ws.Cells.Item(19,1).Value = "'Jun 2026"
$ws.Cells.Item(19,2).Value = "Q3 Rent"
$ws.Cells.Item(19,3).Value = 32097
$ws.Cells.Item(19,4).Value = 6419.4
$ws.Cells.Item(19,5).Value = 38516.4
$ws.Cells.Item(19,6).Value = "Yes"
$ws.Cells.Item(19,7).Value = 248468.6
$ws.Cells.Item(19,8).Value = "Quarter 3 rent payment"

$ws.Cells.Item(20,1).Value = "'Jun 2026"
$ws.Cells.Item(20,2).Value = "Service Charge (Q3)"
$ws.Cells.Item(20,3).Value = 12000
$ws.Cells.Item(20,4).Value = 0
$ws.Cells.Item(20,5).Value = 12000
$ws.Cells.Item(20,6).Value = "No"
$ws.Cells.Item(20,7).Value = 236468.6
$ws.Cells.Item(20,8).Value = "Quarterly service charge"
